$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Cells.Item(5, 6).Value = 329
$ws.Cells.Item(6, 6).Value = 401
$ws.Cells.Item(7, 6).Value = 876
$ws.Cells.Item(9, 6).Value = 523
$ws.Cells.Item(12, 6).Value = 1155
$ws.Cells.Item(15, 6).Value = 37
$ws.Cells.Item(16, 6).Value = 419
$ws.Cells.Item(17, 6).Value = 6677
$ws.Cells.Item(21, 6).Value = 7602
$ws.Cells.Item(22, 6).Value = 41
$ws.Cells.Item(23, 6).Value = 38
$ws.Cells.Item(24, 6).Value = 3404
$ws.Cells.Item(25, 6).Value = 31
$ws.Cells.Item(26, 6).Value = 2118
$ws.Cells.Item(27, 6).Value = 899
$ws.Cells.Item(29, 6).Value = 154
$ws.Cells.Item(30, 6).Value = 353
$ws.Cells.Item(32, 6).Value = 231
$ws.Cells.Item(34, 6).Value = 1724
$ws.Cells.Item(36, 6).Value = 180
$ws.Cells.Item(38, 6).Value = 18
$ws.Cells.Item(39, 6).Value = 1222
$ws.Cells.Item(40, 6).Value = 1818
$ws.Cells.Item(41, 6).Value = 2144
$ws.Cells.Item(42, 6).Value = 10

$ws = $wb.Worksheets.Item("演出")
$ws.Cells.Item(3, 6).Value = 70

$ws = $wb.Worksheets.Item("本地生活")
$ws.Cells.Item(3, 6).Value = 1233

$ws = $wb.Worksheets.Item("全部类型")
$ws.Cells.Item(4, 6).Value = 1233
$ws.Cells.Item(7, 6).Value = 329
$ws.Cells.Item(8, 6).Value = 401
$ws.Cells.Item(9, 6).Value = 876
$ws.Cells.Item(11, 6).Value = 523
$ws.Cells.Item(14, 6).Value = 1155
$ws.Cells.Item(15, 6).Value = 70
$ws.Cells.Item(18, 6).Value = 37
$ws.Cells.Item(19, 6).Value = 419
$ws.Cells.Item(20, 6).Value = 6677
$ws.Cells.Item(24, 6).Value = 7602
$ws.Cells.Item(25, 6).Value = 41
$ws.Cells.Item(26, 6).Value = 38
$ws.Cells.Item(27, 6).Value = 3404
$ws.Cells.Item(28, 6).Value = 31
$ws.Cells.Item(29, 6).Value = 2118
$ws.Cells.Item(30, 6).Value = 899
$ws.Cells.Item(32, 6).Value = 154
$ws.Cells.Item(33, 6).Value = 353
$ws.Cells.Item(36, 6).Value = 232
$ws.Cells.Item(38, 6).Value = 1724
$ws.Cells.Item(40, 6).Value = 180
$ws.Cells.Item(42, 6).Value = 18
$ws.Cells.Item(44, 6).Value = 1222
$ws.Cells.Item(45, 6).Value = 1818
$ws.Cells.Item(47, 6).Value = 2144
$ws.Cells.Item(48, 6).Value = 10
